$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-18 Monday" "2024-03-19 Tuesday"

Replace-Text "719×2=" "946×8="
Replace-Text "351×7=" "845×4="
Replace-Text "395×8=" "706×9="
Replace-Text "220×4=" "955×5="
Replace-Text "387×6=" "279×8="
Replace-Text "999×8=" "253×8="
Replace-Text "617×8=" "367×4="
Replace-Text "830×3=" "323×5="
Replace-Text "660×5=" "935×8="
Replace-Text "662×7=" "801×9="
Replace-Text "141×2=" "936×5="
Replace-Text "330×3=" "692×3="
Replace-Text "969×3=" "892×7="
Replace-Text "985×3=" "241×8="
Replace-Text "223×4=" "977×4="
Replace-Text "361×8=" "446×8="
Replace-Text "604×8=" "217×2="
Replace-Text "216×6=" "332×3="
Replace-Text "794×8=" "370×8="
Replace-Text "268×6=" "376×4="
Replace-Text "168×2=" "422×9="
Replace-Text "259×3=" "931×6="
Replace-Text "353×8=" "150×9="
Replace-Text "393×7=" "493×3="
Replace-Text "592×4=" "839×7="
